$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Add the new "metadata" worksheet right after the existing "data" sheet
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Re-use the bold / bordered / centered+top formatting from the "data" sheet's header row
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$meta.Range("B2").Value = "Additional findings health related - CNV analysis adults"
$meta.Range("C2").Value = 934

$meta.Range("E2").Value = "2021-04-07T10:24:28.854764Z"
$meta.Range("F2").Value = "2021-10-05 14:19:02.013477"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/934/?format=json"

# "data_version" must stay a text value ("2.0"), not be coerced into the number 2.
# Force text via a temporary "@" format, then restore plain/default formatting by
# re-pasting the (unstyled) formats of an untouched cell over it.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.0"
$data.Range("H1").Copy()
$meta.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# Update the "time_taken" column on the "data" sheet to reflect the refreshed run
$data.Range("F2").Value = "2021-10-05 14:19:02.017444"
$data.Range("F3").Value = "2021-10-05 14:19:02.017453"
$data.Range("F4").Value = "2021-10-05 14:19:02.017457"
$data.Range("F5").Value = "2021-10-05 14:19:02.017460"
$data.Range("F6").Value = "2021-10-05 14:19:02.017463"
$data.Range("F7").Value = "2021-10-05 14:19:02.017466"
$data.Range("F8").Value = "2021-10-05 14:19:02.017469"
$data.Range("F9").Value = "2021-10-05 14:19:02.017472"
$data.Range("F10").Value = "2021-10-05 14:19:02.017475"
$data.Range("F11").Value = "2021-10-05 14:19:02.017477"

# Keep "data" as the active/selected sheet (adding "metadata" shouldn't change focus)
$data.Activate() | Out-Null
$data.Range("A1").Select() | Out-Null
